# Add a new student (Allison Serna) to "estudiantes" and a matching
# advisory-session record to "asesorias", enforcing the "one registro per
# day per week" control described in the commit message.

$wb = $excel.ActiveWorkbook

# --- estudiantes sheet: append new student row ---
$wsEstudiantes = $wb.Worksheets.Item("estudiantes")
$wsEstudiantes.Range("A3:K3").NumberFormat = "@"
$wsEstudiantes.Range("A3").Value = "Allison"
$wsEstudiantes.Range("B3").Value = "Serna"
$wsEstudiantes.Range("C3").Value = "1023629824"
$wsEstudiantes.Range("D3").Value = "allis23"
$wsEstudiantes.Range("E3").Value = "allisonserna@gmail.com"
$wsEstudiantes.Range("F3").Value = "3044463002"
$wsEstudiantes.Range("G3").Value = "allison"
$wsEstudiantes.Range("H3").Value = "Estudiante"
$wsEstudiantes.Range("I3").Value = "Antioquia"
$wsEstudiantes.Range("J3").Value = "Medellín"
$wsEstudiantes.Range("K3").Value = "P.C.J.I.C"

# --- asesorias sheet: append new advisory registration row ---
$wsAsesorias = $wb.Worksheets.Item("asesorias")
$wsAsesorias.Range("A2").Value = "Allison Serna"
$wsAsesorias.Range("B2").Value = "allis23"
$wsAsesorias.Range("C2").Value = "s"
$wsAsesorias.Range("D2").Value = "Maryem Ruiz"
$wsAsesorias.Range("E2").Value = "Asesoría PPI"
$wsAsesorias.Range("F2").Value = "23-11-2023"
$wsAsesorias.Range("G2").Value = "02:40 - 03:00"

$wb.Save()
